# Generate Report for Handoff
# Replaces the old GUID-named source file ("73bbbd74-...") with the newly
# generated one ("f046273b-...") throughout the workbook, and bumps the
# associated handoff/handback timestamps.

$newGuid = "f046273b-eaf2-4f5e-bb45-fbe7658c530d"

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-12 05:00:52"

foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newGuid.md"
}

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.8807652c59e58ec4b71b671a5e306c5f743b2ab0.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-12 05:00:46"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

# --- de-de sheet ------------------------------------------------------
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.8807652c59e58ec4b71b671a5e306c5f743b2ab0.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-12 05:00:52"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}
